# Apply updated odds values from FlashScore weekly export
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 7).Value = 1.85  # G2
$ws.Cells.Item(2, 9).Value = 4.1  # I2
$ws.Cells.Item(2, 15).Value = 2.25  # O2
$ws.Cells.Item(2, 16).Value = 2.02  # P2
$ws.Cells.Item(2, 17).Value = 1.62  # Q2
$ws.Cells.Item(2, 18).Value = 5.7  # R2
$ws.Cells.Item(2, 19).Value = 7.7  # S2
$ws.Cells.Item(2, 21).Value = 15  # U2
$ws.Cells.Item(2, 22).Value = 17  # V2
$ws.Cells.Item(2, 24).Value = 7.3  # X2
$ws.Cells.Item(2, 26).Value = 19  # Z2
$ws.Cells.Item(2, 27).Value = 120  # AA2
$ws.Cells.Item(2, 28).Value = 9.25  # AB2
$ws.Cells.Item(2, 29).Value = 21  # AC2
$ws.Cells.Item(2, 30).Value = 14.5  # AD2
$ws.Cells.Item(2, 31).Value = 70  # AE2
$ws.Cells.Item(2, 32).Value = 50  # AF2
$ws.Cells.Item(2, 33).Value = 65  # AG2
# Row 3
$ws.Cells.Item(3, 7).Value = 2.6  # G3
$ws.Cells.Item(3, 9).Value = 2.7  # I3
$ws.Cells.Item(3, 10).Value = 1.33  # J3
$ws.Cells.Item(3, 12).Value = 2.1  # L3
$ws.Cells.Item(3, 13).Value = 1.73  # M3
$ws.Cells.Item(3, 18).Value = 7.5  # R3
$ws.Cells.Item(3, 29).Value = 13  # AC3
$ws.Cells.Item(3, 35).Value = 1.05  # AI3
# Row 4
$ws.Cells.Item(4, 10).Value = 1.33  # J4
$ws.Cells.Item(4, 35).Value = 1.05  # AI4
# Row 5
$ws.Cells.Item(5, 7).Value = 4.2  # G5
$ws.Cells.Item(5, 9).Value = 1.85  # I5
$ws.Cells.Item(5, 10).Value = 1.53  # J5
$ws.Cells.Item(5, 11).Value = 2.32  # K5
$ws.Cells.Item(5, 14).Value = 1.62  # N5
$ws.Cells.Item(5, 15).Value = 2.2  # O5
$ws.Cells.Item(5, 19).Value = 21  # S5
$ws.Cells.Item(5, 35).Value = 1.08  # AI5
# Row 6
$ws.Cells.Item(6, 10).Value = 1.25  # J6
$ws.Cells.Item(6, 12).Value = 1.93  # L6
$ws.Cells.Item(6, 13).Value = 1.93  # M6
$ws.Cells.Item(6, 35).Value = 1.03  # AI6
# Row 7
$ws.Cells.Item(7, 7).Value = 2.25  # G7
$ws.Cells.Item(7, 8).Value = 3.4  # H7
$ws.Cells.Item(7, 9).Value = 3  # I7
$ws.Cells.Item(7, 10).Value = 1.41  # J7
$ws.Cells.Item(7, 11).Value = 2.7  # K7
$ws.Cells.Item(7, 14).Value = 1.53  # N7
$ws.Cells.Item(7, 15).Value = 2.38  # O7
$ws.Cells.Item(7, 19).Value = 10  # S7
$ws.Cells.Item(7, 29).Value = 13  # AC7
$ws.Cells.Item(7, 30).Value = 11  # AD7
$ws.Cells.Item(7, 35).Value = 1.05  # AI7
# Row 8
$ws.Cells.Item(8, 10).Value = 1.41  # J8
$ws.Cells.Item(8, 11).Value = 2.7  # K8
$ws.Cells.Item(8, 35).Value = 1.05  # AI8
# Row 9
$ws.Cells.Item(9, 35).Value = 1.08  # AI9
$ws.Cells.Item(9, 36).Value = 8  # AJ9
# Row 11
$ws.Cells.Item(11, 8).Value = 3.75  # H11
$ws.Cells.Item(11, 10).Value = 1.29  # J11
$ws.Cells.Item(11, 11).Value = 3.75  # K11
# Row 12
$ws.Cells.Item(12, 36).Value = 13  # AJ12
# Row 13
$ws.Cells.Item(13, 7).Value = 2.75  # G13
$ws.Cells.Item(13, 8).Value = 3.1  # H13
$ws.Cells.Item(13, 10).Value = 1.5  # J13
$ws.Cells.Item(13, 11).Value = 2.5  # K13
$ws.Cells.Item(13, 12).Value = 2.5  # L13
$ws.Cells.Item(13, 13).Value = 1.5  # M13
$ws.Cells.Item(13, 14).Value = 1.57  # N13
$ws.Cells.Item(13, 15).Value = 2.25  # O13
$ws.Cells.Item(13, 16).Value = 2.05  # P13
$ws.Cells.Item(13, 17).Value = 1.7  # Q13
$ws.Cells.Item(13, 24).Value = 7  # X13
$ws.Cells.Item(13, 26).Value = 19  # Z13
$ws.Cells.Item(13, 28).Value = 6.5  # AB13
$ws.Cells.Item(13, 29).Value = 11  # AC13
$ws.Cells.Item(13, 32).Value = 26  # AF13
$ws.Cells.Item(13, 34).Value = 501  # AH13
$ws.Cells.Item(13, 35).Value = 1.1  # AI13
$ws.Cells.Item(13, 36).Value = 7  # AJ13
# Row 14
$ws.Cells.Item(14, 9).Value = 2.55  # I14
$ws.Cells.Item(14, 10).Value = 1.36  # J14
$ws.Cells.Item(14, 11).Value = 3  # K14
$ws.Cells.Item(14, 12).Value = 2.15  # L14
$ws.Cells.Item(14, 13).Value = 1.67  # M14
$ws.Cells.Item(14, 18).Value = 8  # R14
$ws.Cells.Item(14, 19).Value = 13  # S14
$ws.Cells.Item(14, 24).Value = 8.5  # X14
$ws.Cells.Item(14, 35).Value = 1.07  # AI14
$ws.Cells.Item(14, 36).Value = 9  # AJ14
# Row 15
$ws.Cells.Item(15, 10).Value = 1.3  # J15
$ws.Cells.Item(15, 11).Value = 3.4  # K15
$ws.Cells.Item(15, 12).Value = 1.97  # L15
$ws.Cells.Item(15, 13).Value = 1.77  # M15
$ws.Cells.Item(15, 24).Value = 10  # X15
$ws.Cells.Item(15, 30).Value = 10  # AD15
$ws.Cells.Item(15, 35).Value = 1.05  # AI15
$ws.Cells.Item(15, 36).Value = 11  # AJ15
# Row 17
$ws.Cells.Item(17, 7).Value = 3.7  # G17
$ws.Cells.Item(17, 9).Value = 2.05  # I17
$ws.Cells.Item(17, 16).Value = 1.91  # P17
$ws.Cells.Item(17, 17).Value = 1.91  # Q17
$ws.Cells.Item(17, 18).Value = 10  # R17
$ws.Cells.Item(17, 19).Value = 19  # S17
$ws.Cells.Item(17, 20).Value = 13  # T17
$ws.Cells.Item(17, 21).Value = 41  # U17
$ws.Cells.Item(17, 23).Value = 41  # W17
$ws.Cells.Item(17, 28).Value = 7  # AB17
$ws.Cells.Item(17, 29).Value = 9.5  # AC17
$ws.Cells.Item(17, 31).Value = 17  # AE17
$ws.Cells.Item(17, 32).Value = 17  # AF17
$ws.Cells.Item(17, 35).Value = 1.07  # AI17
$ws.Cells.Item(17, 36).Value = 9  # AJ17
# Row 18
$ws.Cells.Item(18, 7).Value = 1.75  # G18
$ws.Cells.Item(18, 8).Value = 3.5  # H18
$ws.Cells.Item(18, 9).Value = 4.33  # I18
$ws.Cells.Item(18, 10).Value = 1.36  # J18
$ws.Cells.Item(18, 11).Value = 3  # K18
$ws.Cells.Item(18, 16).Value = 2  # P18
$ws.Cells.Item(18, 17).Value = 1.73  # Q18
$ws.Cells.Item(18, 18).Value = 6  # R18
$ws.Cells.Item(18, 19).Value = 8  # S18
$ws.Cells.Item(18, 21).Value = 15  # U18
$ws.Cells.Item(18, 26).Value = 19  # Z18
$ws.Cells.Item(18, 27).Value = 67  # AA18
$ws.Cells.Item(18, 28).Value = 11  # AB18
$ws.Cells.Item(18, 29).Value = 21  # AC18
$ws.Cells.Item(18, 34).Value = 401  # AH18
# Row 20
$ws.Cells.Item(20, 8).Value = 6.8  # H20
$ws.Cells.Item(20, 9).Value = 16  # I20
$ws.Cells.Item(20, 16).Value = 2.14  # P20
$ws.Cells.Item(20, 17).Value = 1.63  # Q20
$ws.Cells.Item(20, 18).Value = 9  # R20
$ws.Cells.Item(20, 19).Value = 6.1  # S20
$ws.Cells.Item(20, 20).Value = 9.5  # T20
$ws.Cells.Item(20, 21).Value = 5.8  # U20
$ws.Cells.Item(20, 23).Value = 25  # W20
$ws.Cells.Item(20, 24).Value = 21  # X20
$ws.Cells.Item(20, 25).Value = 14.5  # Y20
$ws.Cells.Item(20, 27).Value = 90  # AA20
$ws.Cells.Item(20, 28).Value = 45  # AB20
$ws.Cells.Item(20, 29).Value = 150  # AC20
$ws.Cells.Item(20, 30).Value = 45  # AD20
$ws.Cells.Item(20, 32).Value = 200  # AF20
$ws.Cells.Item(20, 33).Value = 110  # AG20
# Row 21
$ws.Cells.Item(21, 7).Value = 5.2  # G21
$ws.Cells.Item(21, 8).Value = 3.9  # H21
$ws.Cells.Item(21, 9).Value = 1.52  # I21
$ws.Cells.Item(21, 11).Value = 4.4  # K21
$ws.Cells.Item(21, 13).Value = 2.2  # M21
$ws.Cells.Item(21, 16).Value = 1.66  # P21
$ws.Cells.Item(21, 17).Value = 2.09  # Q21
$ws.Cells.Item(21, 18).Value = 15  # R21
$ws.Cells.Item(21, 19).Value = 29  # S21
$ws.Cells.Item(21, 20).Value = 13.5  # T21
$ws.Cells.Item(21, 21).Value = 75  # U21
$ws.Cells.Item(21, 22).Value = 37  # V21
$ws.Cells.Item(21, 23).Value = 32  # W21
$ws.Cells.Item(21, 24).Value = 14  # X21
$ws.Cells.Item(21, 25).Value = 7  # Y21
$ws.Cells.Item(21, 26).Value = 11.75  # Z21
$ws.Cells.Item(21, 28).Value = 7.5  # AB21
$ws.Cells.Item(21, 29).Value = 7.3  # AC21
$ws.Cells.Item(21, 30).Value = 6.8  # AD21
$ws.Cells.Item(21, 31).Value = 10  # AE21
$ws.Cells.Item(21, 33).Value = 16  # AG21
$ws.Cells.Item(21, 34).Value = 200  # AH21
# Row 23
$ws.Cells.Item(23, 7).Value = 3.6  # G23
$ws.Cells.Item(23, 9).Value = 1.95  # I23
$ws.Cells.Item(23, 12).Value = 1.9  # L23
$ws.Cells.Item(23, 16).Value = 1.73  # P23
$ws.Cells.Item(23, 20).Value = 13  # T23
$ws.Cells.Item(23, 29).Value = 9.5  # AC23
$ws.Cells.Item(23, 31).Value = 17  # AE23
$ws.Cells.Item(23, 32).Value = 15  # AF23
# Row 24
$ws.Cells.Item(24, 16).Value = 1.83  # P24
$ws.Cells.Item(24, 17).Value = 1.83  # Q24
# Row 25
$ws.Cells.Item(25, 12).Value = 1.88  # L25
$ws.Cells.Item(25, 13).Value = 1.98  # M25
$ws.Cells.Item(25, 16).Value = 1.73  # P25
# Row 26
$ws.Cells.Item(26, 13).Value = 2  # M26
$ws.Cells.Item(26, 16).Value = 1.67  # P26
# Row 27
$ws.Cells.Item(27, 8).Value = 7  # H27
$ws.Cells.Item(27, 9).Value = 23  # I27
$ws.Cells.Item(27, 16).Value = 2.32  # P27
$ws.Cells.Item(27, 17).Value = 1.54  # Q27
$ws.Cells.Item(27, 25).Value = 16  # Y27
$ws.Cells.Item(27, 28).Value = 70  # AB27
$ws.Cells.Item(27, 30).Value = 70  # AD27
# Row 29
$ws.Cells.Item(29, 7).Value = 1.2  # G29
$ws.Cells.Item(29, 8).Value = 5.8  # H29
$ws.Cells.Item(29, 9).Value = 14  # I29
$ws.Cells.Item(29, 10).Value = 1.28  # J29
$ws.Cells.Item(29, 11).Value = 3.45  # K29
$ws.Cells.Item(29, 12).Value = 1.86  # L29
$ws.Cells.Item(29, 13).Value = 1.88  # M29
$ws.Cells.Item(29, 15).Value = 2.8  # O29
$ws.Cells.Item(29, 16).Value = 2.9  # P29
$ws.Cells.Item(29, 17).Value = 1.38  # Q29
$ws.Cells.Item(29, 19).Value = 3.75  # S29
$ws.Cells.Item(29, 20).Value = 8.199999999999999  # T29
$ws.Cells.Item(29, 21).Value = 4.9  # U29
$ws.Cells.Item(29, 23).Value = 50  # W29
$ws.Cells.Item(29, 24).Value = 9  # X29
$ws.Cells.Item(29, 25).Value = 9.800000000000001  # Y29
$ws.Cells.Item(29, 28).Value = 20  # AB29
$ws.Cells.Item(29, 29).Value = 90  # AC29
$ws.Cells.Item(29, 30).Value = 40  # AD29
# Row 30
$ws.Cells.Item(30, 7).Value = 2.15  # G30
$ws.Cells.Item(30, 9).Value = 3.1  # I30
$ws.Cells.Item(30, 18).Value = 7  # R30
$ws.Cells.Item(30, 32).Value = 29  # AF30
$ws.Cells.Item(30, 33).Value = 41  # AG30
$ws.Cells.Item(30, 35).Value = 1.06  # AI30
$ws.Cells.Item(30, 36).Value = 10  # AJ30
# Row 31
$ws.Cells.Item(31, 10).Value = 1.25  # J31
$ws.Cells.Item(31, 11).Value = 3.75  # K31
$ws.Cells.Item(31, 12).Value = 1.85  # L31
$ws.Cells.Item(31, 13).Value = 2  # M31
# Row 32
$ws.Cells.Item(32, 7).Value = 1.42  # G32
$ws.Cells.Item(32, 19).Value = 6  # S32
$ws.Cells.Item(32, 20).Value = 9  # T32
$ws.Cells.Item(32, 25).Value = 8.5  # Y32
$ws.Cells.Item(32, 27).Value = 81  # AA32
$ws.Cells.Item(32, 35).Value = 1.05  # AI32
$ws.Cells.Item(32, 36).Value = 11  # AJ32
# Row 33
$ws.Cells.Item(33, 10).Value = 1.4  # J33
$ws.Cells.Item(33, 11).Value = 2.75  # K33
$ws.Cells.Item(33, 12).Value = 2.15  # L33
$ws.Cells.Item(33, 13).Value = 1.62  # M33
$ws.Cells.Item(33, 14).Value = 1.44  # N33
$ws.Cells.Item(33, 15).Value = 2.62  # O33
$ws.Cells.Item(33, 16).Value = 1.83  # P33
$ws.Cells.Item(33, 17).Value = 1.87  # Q33
$ws.Cells.Item(33, 18).Value = 7.5  # R33
$ws.Cells.Item(33, 19).Value = 13  # S33
$ws.Cells.Item(33, 20).Value = 10  # T33
$ws.Cells.Item(33, 22).Value = 25  # V33
$ws.Cells.Item(33, 23).Value = 37  # W33
$ws.Cells.Item(33, 24).Value = 6.2  # X33
$ws.Cells.Item(33, 26).Value = 14.5  # Z33
$ws.Cells.Item(33, 27).Value = 75  # AA33
$ws.Cells.Item(33, 28).Value = 7.6  # AB33
$ws.Cells.Item(33, 29).Value = 13  # AC33
$ws.Cells.Item(33, 30).Value = 9.75  # AD33
$ws.Cells.Item(33, 32).Value = 24  # AF33
$ws.Cells.Item(33, 34).Value = 700  # AH33
$ws.Cells.Item(33, 36).Value = 6.2  # AJ33
